$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Title paragraph: merge the split "B" / "lueprint Summary for {{" runs
#    into a single run "Blueprint Summary for {{" and drop the stray
#    _GoBack bookmark that used to sit between them.
# ---------------------------------------------------------------------------
$rTitle = $d.Content
$rTitle.Find.Execute("Blueprint Summary for {{", $false, $false, $false, $false, $false, $true, 1, $false, "Blueprint Summary for {{", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) {{SLR}} -> {{table:slr}}, split across 3 runs ({{ / table:slr / }}),
#    with a _GoBack bookmark placed between the 2nd and 3rd run.
# ---------------------------------------------------------------------------
$rSlr = $d.Content
$foundSlr = $rSlr.Find.Execute("{{SLR}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundSlr) {
    throw "could not find {{SLR}} placeholder"
}
$slrStart = $rSlr.Start
$rSlr.Text = "{{table:slr}}"

$slrMid = $d.Range($slrStart + 2, $slrStart + 2 + 9)   # "table:slr"
$slrMid.Font.Color = 255
$slrMid.Font.Color = 3355443

$slrBmPos = $slrStart + 2 + 9
$slrBmRange = $d.Range($slrBmPos, $slrBmPos)
$d.Bookmarks.Add("_GoBack", $slrBmRange) | Out-Null

# ---------------------------------------------------------------------------
# 3) {{URBAN}} -> {{table:urban}}, split across 3 runs
#    ({{ / table:urban / }}).
# ---------------------------------------------------------------------------
$rUrban = $d.Content
$foundUrban = $rUrban.Find.Execute("{{URBAN}}", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundUrban) {
    throw "could not find {{URBAN}} placeholder"
}
$urbanStart = $rUrban.Start
$rUrban.Text = "{{table:urban}}"

$urbanMid = $d.Range($urbanStart + 2, $urbanStart + 2 + 11)   # "table:urban"
$urbanMid.Font.Color = 255
$urbanMid.Font.Color = 3355443

Write-Output "done"
